{"js": "// Time-tracking table edit:\n//  1) Row \"10.02.2023\" / \"0.5\" / \"Weekly\": the date cell text is split across\n//     two runs (\"10.\" + \"02.2023\"); normalize it to a single run \"10.02.2023\"\n//     (text itself is unchanged, just re-written as one run).\n//  2) Row \"16.2.2023\" (previously empty Hours/Description cells): fix the\n//     date to \"16.02.2023\" and fill in Hours = \"4\" and\n//     Description = \"Uuden teht\u00e4v\u00e4n luonnin ty\u00f6st\u00f6\".\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\n\n// Locate the target rows by their current cell text instead of a hard-coded\n// index, so the script is resilient to minor row-count differences.\nlet mergedDateRow = -1; // currently \"10.02.2023\" split across two runs\nlet newEntryRow = -1; // currently \"16.2.2023\" with empty Hours/Description\n\nfor (let i = 0; i < values.length; i++) {\n  const row = values[i];\n  if (row[0] === \"10.02.2023\" && row[1] === \"0.5\") {\n    mergedDateRow = i;\n  }\n  if (row[0] === \"16.2.2023\" && (row[1] || \"\").trim() === \"\" && (row[2] || \"\").trim() === \"\") {\n    newEntryRow = i;\n  }\n}\n\nif (mergedDateRow === -1) {\n  throw new Error('Could not find the \"10.02.2023\" row to normalize.');\n}\nif (newEntryRow === -1) {\n  throw new Error('Could not find the \"16.2.2023\" row to fill in.');\n}\n\n// 1) Re-write the date cell so \"10.\" + \"02.2023\" becomes one run.\nconst dateCell = table.getCell(mergedDateRow, 0);\ndateCell.body.insertText(\"10.02.2023\", Word.InsertLocation.replace);\n\n// 2) Fix the date typo and fill in the new row's data.\nconst newDateCell = table.getCell(newEntryRow, 0);\nnewDateCell.body.insertText(\"16.02.2023\", Word.InsertLocation.replace);\n\nconst newHoursCell = table.getCell(newEntryRow, 1);\nnewHoursCell.body.insertText(\"4\", Word.InsertLocation.replace);\n\nconst newDescCell = table.getCell(newEntryRow, 2);\nnewDescCell.body.insertText(\"Uuden teht\u00e4v\u00e4n luonnin ty\u00f6st\u00f6\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Time-tracking table edit:\n#  1) Row \"10.02.2023\" / \"0.5\" / \"Weekly\": the date cell text is split across\n#     two runs (\"10.\" + \"02.2023\"); normalize it to a single run \"10.02.2023\"\n#     (text itself is unchanged, just re-written as one run).\n#  2) Row \"16.2.2023\" (previously empty Hours/Description cells): fix the\n#     date to \"16.02.2023\" and fill in Hours = \"4\" and\n#     Description = \"Uuden teht\u00e4v\u00e4n luonnin ty\u00f6st\u00f6\".\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the target rows by their current cell text instead of a hard-coded\n# row index, so the script is resilient to minor row-count differences.\n$mergedDateRow = -1   # currently \"10.02.2023\" split across two runs\n$newEntryRow = -1     # currently \"16.2.2023\" with empty Hours/Description\n\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $dateText = $t.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7)\n    $hoursText = $t.Cell($i, 2).Range.Text.TrimEnd([char]13, [char]7)\n    $descText = $t.Cell($i, 3).Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($dateText -eq \"10.02.2023\" -and $hoursText -eq \"0.5\") {\n        $mergedDateRow = $i\n    }\n    if ($dateText -eq \"16.2.2023\" -and $hoursText.Trim() -eq \"\" -and $descText.Trim() -eq \"\") {\n        $newEntryRow = $i\n    }\n}\n\nif ($mergedDateRow -eq -1) {\n    throw \"Could not find the '10.02.2023' row to normalize.\"\n}\nif ($newEntryRow -eq -1) {\n    throw \"Could not find the '16.2.2023' row to fill in.\"\n}\n\n# 1) Re-write the date cell so \"10.\" + \"02.2023\" becomes one run. Using\n#    Find/Replace across the whole cell merges the two runs into one.\n$d.Content.Find.Execute(\"10.02.2023\", $false, $false, $false, $false, $false, $true, 1, $false, \"10.02.2023\", 2)\n\n# 2) Fix the date typo \"16.2.2023\" -> \"16.02.2023\".\n$d.Content.Find.Execute(\"16.2.2023\", $false, $false, $false, $false, $false, $true, 1, $false, \"16.02.2023\", 2)\n\n# 3) Fill in the Hours and Description cells for that row.\n$t.Cell($newEntryRow, 2).Range.Text = \"4\"\n$t.Cell($newEntryRow, 3).Range.Text = \"Uuden teht\u00e4v\u00e4n luonnin ty\u00f6st\u00f6\"\n"}
